$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -----------------------------------------------------------------
# Insert one row at row 27 - this pushes the existing "SHIMS" block
# (old rows 28-34) down by one row (new rows 29-35), exactly like the
# target diff, while Excel naturally carries the per-cell formatting
# along with the shifted rows. Do this FIRST so the row numbers used
# below for the new content (37-47) are not themselves shifted.
# -----------------------------------------------------------------
$ws.Rows("27:27").Insert()

# -----------------------------------------------------------------
# New "SCOPE OF SHIMS" section (rows 37-41)
# -----------------------------------------------------------------
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuse the highlighted style

$ws.Range("B37").Value = "SCOPE OF SHIMS"
$ws.Range("B38").Value = "Defined shims are used only inside the SHIMS context i.e. with in the below USING block"
$ws.Range("C39").Value = "using(ShimsContext.Create())"
$ws.Range("C40").Value = "            {"
$ws.Range("C41").Value = "   }"

# New note that lands in row 27
$ws.Range("C27").Value = "For STUBS, we need interfaces and so we have to inject into the classess"

# -----------------------------------------------------------------
# New "ShimPatientService" sample section (rows 44-47)
# -----------------------------------------------------------------
$ws.Range("B44").Value = 'WHAT IS THE PURPOSE OF "ShimPatientService"   VS "ShimPatientService.AllInstances"?'
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null

$ws.Range("A21").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null

$ws.Range("B47").Value = "ShimPatientService.AllInstances: "
$ws.Range("B46").Value = "ShimPatientService:    "

$ws.Range("F47").Value = "This is to shim(provide alternative implementation) for the INSTANCE & PRIVATE Methods/Properties"
$ws.Range("F46").Value = "This is to shim(provide alternative implementation) for the STATIC Methods/Properties"

$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# Update the sheet view to match the final scroll/selection state
# -----------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("F47").Select()
